$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "Fig_3a"
$wb.Worksheets.Item(2).Name = "Fig_3b"
$wb.Worksheets.Item(3).Name = "Fig_3c"
$wb.Worksheets.Item(4).Name = "Fig_3d"
$wb.Worksheets.Item(5).Name = "Fig_3e"

$wb.Worksheets.Item(5).Activate()
